# Natmi following Dr Hou advice
# Rebuild the Hgf-Met LR-pair table: add the missing "ECs" sending/target
# cluster so the sheet now covers the full 3x3 cross of {ECs, FAPs, sCs}
# (rows 2-10), replacing the old 2x2 data that only had {FAPs, sCs}.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster | E-T metrics

# row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hgf"
$ws.Range("C2").Value = "Met"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 5.928568666666667
$ws.Range("H2").Value = 17.785706
$ws.Range("I2").Value = 0.3809768389628236
$ws.Range("J2").Value = 0.3809768389628236
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.254050666666667
$ws.Range("N2").Value = 6.762152
$ws.Range("O2").Value = 0.5759961839619929
$ws.Range("P2").Value = 0.575996183961993
$ws.Range("Q2").Value = 13.36329415547911
$ws.Range("R2").Value = 120.269647399312
$ws.Range("S2").Value = 0.2194412054204891
$ws.Range("T2").Value = 0.2194412054204892

# row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hgf"
$ws.Range("C3").Value = "Met"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 5.928568666666667
$ws.Range("H3").Value = 17.785706
$ws.Range("I3").Value = 0.3809768389628236
$ws.Range("J3").Value = 0.3809768389628236
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05371366666666667
$ws.Range("N3").Value = 0.161141
$ws.Range("O3").Value = 0.01372589688605336
$ws.Range("P3").Value = 0.01372589688605336
$ws.Range("Q3").Value = 0.3184451611717778
$ws.Range("R3").Value = 2.866006450546
$ws.Range("S3").Value = 0.005229248807578274
$ws.Range("T3").Value = 0.005229248807578274

# row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hgf"
$ws.Range("C4").Value = "Met"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 5.928568666666667
$ws.Range("H4").Value = 17.785706
$ws.Range("I4").Value = 0.3809768389628236
$ws.Range("J4").Value = 0.3809768389628236
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.605544
$ws.Range("N4").Value = 4.816632
$ws.Range("O4").Value = 0.4102779191519537
$ws.Range("P4").Value = 0.4102779191519537
$ws.Range("Q4").Value = 9.518577851354667
$ws.Range("R4").Value = 85.66720066219202
$ws.Range("S4").Value = 0.1563063847347562
$ws.Range("T4").Value = 0.1563063847347562

# row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hgf"
$ws.Range("C5").Value = "Met"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.583521
$ws.Range("H5").Value = 25.750563
$ws.Range("I5").Value = 0.5515872180307627
$ws.Range("J5").Value = 0.5515872180307626
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.254050666666667
$ws.Range("N5").Value = 6.762152
$ws.Range("O5").Value = 0.5759961839619929
$ws.Range("P5").Value = 0.575996183961993
$ws.Range("Q5").Value = 19.34769123239733
$ws.Range("R5").Value = 174.129221091576
$ws.Range("S5").Value = 0.3177121327079311
$ws.Range("T5").Value = 0.3177121327079311

# row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hgf"
$ws.Range("C6").Value = "Met"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.583521
$ws.Range("H6").Value = 25.750563
$ws.Range("I6").Value = 0.5515872180307627
$ws.Range("J6").Value = 0.5515872180307626
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05371366666666667
$ws.Range("N6").Value = 0.161141
$ws.Range("O6").Value = 0.01372589688605336
$ws.Range("P6").Value = 0.01372589688605336
$ws.Range("Q6").Value = 0.4610523858203333
$ws.Range("R6").Value = 4.149471472383
$ws.Range("S6").Value = 0.007571029278355282
$ws.Range("T6").Value = 0.007571029278355281

# row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hgf"
$ws.Range("C7").Value = "Met"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.583521
$ws.Range("H7").Value = 25.750563
$ws.Range("I7").Value = 0.5515872180307627
$ws.Range("J7").Value = 0.5515872180307626
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.605544
$ws.Range("N7").Value = 4.816632
$ws.Range("O7").Value = 0.4102779191519537
$ws.Range("P7").Value = 0.4102779191519537
$ws.Range("Q7").Value = 13.781220640424
$ws.Range("R7").Value = 124.030985763816
$ws.Range("S7").Value = 0.2263040560444763
$ws.Range("T7").Value = 0.2263040560444763

# row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Hgf"
$ws.Range("C8").Value = "Met"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.049404
$ws.Range("H8").Value = 3.148212
$ws.Range("I8").Value = 0.06743594300641363
$ws.Range("J8").Value = 0.06743594300641362
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.254050666666667
$ws.Range("N8").Value = 6.762152
$ws.Range("O8").Value = 0.5759961839619929
$ws.Range("P8").Value = 0.575996183961993
$ws.Range("Q8").Value = 2.365409785802667
$ws.Range("R8").Value = 21.288688072224
$ws.Range("S8").Value = 0.0388428458335727
$ws.Range("T8").Value = 0.0388428458335727

# row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Hgf"
$ws.Range("C9").Value = "Met"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.049404
$ws.Range("H9").Value = 3.148212
$ws.Range("I9").Value = 0.06743594300641363
$ws.Range("J9").Value = 0.06743594300641362
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05371366666666667
$ws.Range("N9").Value = 0.161141
$ws.Range("O9").Value = 0.01372589688605336
$ws.Range("P9").Value = 0.01372589688605336
$ws.Range("Q9").Value = 0.05636733665466667
$ws.Range("R9").Value = 0.507306029892
$ws.Range("S9").Value = 0.0009256188001198047
$ws.Range("T9").Value = 0.0009256188001198047

# row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Hgf"
$ws.Range("C10").Value = "Met"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.049404
$ws.Range("H10").Value = 3.148212
$ws.Range("I10").Value = 0.06743594300641363
$ws.Range("J10").Value = 0.06743594300641362
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.605544
$ws.Range("N10").Value = 4.816632
$ws.Range("O10").Value = 0.4102779191519537
$ws.Range("P10").Value = 0.4102779191519537
$ws.Range("Q10").Value = 1.684864295776
$ws.Range("R10").Value = 15.163778661984
$ws.Range("S10").Value = 0.02766747837272113
$ws.Range("T10").Value = 0.02766747837272112

